{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Paragraph 0: title/date line - replace text (drops the <w:br/> + second title run)\nparagraphs.items[0].insertText(\"\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 -07.11.24: \u26a1\ufe0f\ud83d\ude80\", \"Replace\");\n\nparagraphs.items[1].insertText(\"Cross-layer Attention Sharing for Large Language Models\", \"Replace\");\nparagraphs.items[2].insertText(\"\u05d0\u05ea\u05dd \u05d1\u05d8\u05d7 \u05d9\u05d5\u05d3\u05e2\u05d9\u05dd \u05d4\u05e8\u05e6\u05d4 \u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4 \u05e2\u05dc\u05d5\u05dc \u05dc\u05d4\u05d9\u05d5\u05ea \u05d3\u05d1\u05e8 \u05d3\u05d9 \u05d9\u05e7\u05e8 \u05de\u05d1\u05d7\u05d9\u05e0\u05ea \u05de\u05e9\u05d0\u05d1\u05d9 \u05d7\u05d9\u05e9\u05d5\u05d1 \u05d5\u05d2\u05dd \u05d4\u05d6\u05db\u05e8\u05d5\u05df. \u05d1\u05d8\u05d7 \u05db\u05d0\u05e9\u05e8 \u05d9\u05e9 \u05dc\u05db\u05dd \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05e2\u05dd \u05e2\u05e9\u05e8\u05d5\u05ea \u05de\u05d9\u05dc\u05d9\u05d0\u05e8\u05d3\u05d9 \u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd \u05e2\u05dc \u05e2\u05e9\u05e8\u05d5\u05ea \u05e8\u05d1\u05d5\u05ea \u05e9\u05dc \u05e9\u05db\u05d1\u05d5\u05ea \u05e9\u05dc \u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8\u05d9\u05dd. \u05d0\u05d7\u05d3 \u05d4\u05d3\u05d1\u05e8\u05d9\u05dd \u05d4\u05db\u05d1\u05d3\u05d9\u05dd \u05e9\u05de\u05e6\u05e8\u05d9\u05db\u05d9\u05dd \u05dc\u05d0 \u05de\u05e2\u05d8 \u05d6\u05d9\u05db\u05e8\u05d5\u05df \u05d4\u05d5\u05d0 KV-Cache, \u05e9\u05d1\u05d5 \u05e0\u05e9\u05de\u05e8\u05d9\u05dd \u05d4\u05de\u05db\u05e4\u05dc\u05d5\u05ea \u05e9\u05dc \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 (\u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2\u05e1) \u05e9\u05dc \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d1\u05de\u05d8\u05e8\u05d9\u05e6\u05d5\u05ea K \u05d5- V \u05dc\u05db\u05dc \u05d4\u05e9\u05db\u05d1\u05d5\u05ea \u05d5\u05dc\u05db\u05dc \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05e9\u05db\u05d1\u05e8 \u05d2\u05d5\u05e0\u05e8\u05d8\u05d5 (\u05db\u05d5\u05dc\u05dc \u05d4\u05e4\u05e8\u05d5\u05de\u05e4\u05d8 - \u05de\u05d3\u05d5\u05d1\u05e8 \u05d1\u05de\u05d5\u05d3\u05dc\u05d9 \u05d4\u05d3\u05e7\u05d5\u05d3\u05e8\u05d9\u05dd).\", \"Replace\");\nparagraphs.items[3].insertText(\"\u05db\u05de\u05d5\u05d1\u05df \u05e9\u05db\u05d0\u05e9\u05e8 \u05d4\u05de\u05d9\u05de\u05d3\u05d9\u05dd \u05e9\u05dc \u05d5\u05e7\u05d8\u05d5\u05e8\u05d9 \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05d5\u05d4\u05de\u05d8\u05e8\u05d9\u05e6\u05d5\u05ea \u05dc\u05d0 \u05e7\u05d8\u05e0\u05d9\u05dd \u05d5\u05d2\u05dd \u05d0\u05d5\u05e8\u05da \u05d4\u05d4\u05e7\u05e9\u05e8 \u05e0\u05de\u05d3\u05d3 \u05d1\u05e2\u05e9\u05e8\u05d5\u05ea \u05d5\u05de\u05d0\u05d5\u05ea \u05d0\u05dc\u05e4\u05d9\u05dd KV-Cache \u05d3\u05d5\u05e8\u05e9 \u05d4\u05e8\u05d1\u05d4 \u05de\u05d0\u05d5\u05d3 \u05d6\u05d9\u05db\u05e8\u05d5\u05df. \u05d1\u05e2\u05d1\u05e8 \u05d9\u05e6\u05d0\u05d5 \u05dc\u05d0 \u05de\u05e2\u05d8 \u05de\u05d0\u05de\u05e8\u05d9\u05dd \u05e9\u05e0\u05d9\u05e1\u05d5 \u05dc\u05d3\u05d7\u05d5\u05e1 \u05d0\u05d5\u05ea\u05d5 \u05e2\u05dc \u05d9\u05d3\u05d9 \u05e0\u05d9\u05ea\u05d5\u05d7 \u05d5\u05d6\u05d9\u05d4\u05d5\u05d9 \u05d9\u05ea\u05d9\u05e8\u05d5\u05d9\u05d5\u05ea \u05d0\u05d1\u05dc \u05d6\u05d4 \u05d1\u05d3\u05f4\u05db \u05e0\u05e2\u05e9\u05d4 \u05e4\u05e8 \u05e9\u05db\u05d1\u05d4 (= \u05d1\u05dc\u05d5\u05e7 \u05d4\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8). \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05de\u05e1\u05d5\u05e7\u05e8 \u05de\u05e6\u05d9\u05e2 \u05dc\u05d4\u05ea\u05d1\u05d5\u05e0\u05df \u05d1\u05d3\u05d7\u05d9\u05e1\u05ea KV-cache \u05de\u05e4\u05e8\u05e1\u05e4\u05e7\u05d8\u05d9\u05d1\u05d4 \u05e8\u05d7\u05d1\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05d5\u05dc\u05e0\u05e1\u05d5\u05ea \u05dc\u05d3\u05d7\u05d5\u05e1 \u05d0\u05d5\u05ea\u05d5 \u05d3\u05e8\u05da \u05e0\u05d9\u05e6\u05d5\u05dc \u05d4\u05ea\u05dc\u05d5\u05d9\u05d5\u05ea \u05e9\u05dc \u05d4-KV-cache \u05d1\u05d9\u05df \u05d4\u05e9\u05db\u05d1\u05d5\u05ea \u05d4\u05e9\u05d5\u05e0\u05d5\u05ea.\", \"Replace\");\nparagraphs.items[4].insertText(\"\u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05d7\u05e7\u05e8\u05d5 \u05d3\u05de\u05d9\u05d5\u05df \u05d1\u05d9\u05df \u05d4\u05d7\u05dc\u05e7\u05d9\u05dd \u05d4\u05e9\u05d5\u05e0\u05d9\u05dd \u05d1\u05d1\u05dc\u05d5\u05e7 \u05d4\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8\u05d9\u05dd  (\u05de\u05db\u05e4\u05dc\u05d5\u05ea \u05e9\u05dc \u05d4\u05de\u05d8\u05e8\u05d9\u05e6\u05d5\u05ea \u05d4\u05e9\u05d5\u05e0\u05d5\u05ea \u05d1\u05d5\u05e7\u05d8\u05d5\u05e8\u05d9 \u05d9\u05d9\u05e6\u05d5\u05d2, \u05de\u05e7\u05d3\u05de\u05d9 attention \u05d5\u05db\u05d3\u05d5\u05de\u05d4) \u05d5\u05d4\u05d2\u05d9\u05e2\u05d5 \u05dc\u05de\u05e1\u05e7\u05e0\u05d4 \u05e9\u05e0\u05d9\u05ea\u05df \u05f4\u05dc\u05d4\u05e1\u05d9\u05e7\u05f4 \u05d0\u05ea \u05de\u05e7\u05d3\u05de\u05d9 \u05d4-attention \u05e9\u05dc \u05e9\u05db\u05d1\u05d4 n \u05de\u05d4\u05d3\u05d0\u05d8\u05d4 \u05e9\u05dc \u05e9\u05db\u05d1\u05d4 n-1 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d7\u05e1\u05db\u05d5\u05e0\u05d9\u05ea \u05d7\u05d9\u05e9\u05d5\u05d1\u05d9\u05ea. \u05db\u05dc\u05d5\u05de\u05e8 \u05e2\u05dd \u05d4\u05e8\u05d1\u05d4 \u05e4\u05d7\u05d5\u05ea \u05de\u05e9\u05e7\u05d5\u05dc\u05d5\u05ea \u05de\u05d4\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05de\u05e8 \u05d4\u05e8\u05d2\u05d9\u05dc. \u05db\u05dc\u05d5\u05de\u05e8 \u05d4\u05d4\u05e6\u05e2\u05d4 \u05d4\u05d9\u05d0 \u05dc\u05e2\u05e9\u05d5\u05ea \u05e1\u05d5\u05d2 \u05e9\u05dc  LoRa \u05d0\u05d1\u05dc \u05dc\u05de\u05e7\u05d3\u05de\u05d9 \u05d4-attention. \", \"Replace\");\nparagraphs.items[5].insertText(\"\u05d1\u05e6\u05d5\u05e8\u05d4 \u05e7\u05e6\u05ea \u05d9\u05d5\u05ea\u05e8 \u05e7\u05d5\u05e0\u05e7\u05e8\u05d8\u05d9\u05ea \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d7\u05dc\u05d9\u05e3 \u05de\u05d8\u05e8\u05d9\u05e6\u05d5\u05ea W_Q \u05d5-W_K \u05d1\u05de\u05d8\u05e8\u05d9\u05e6\u05d5\u05ea \u05d1\u05e2\u05dc\u05d5\u05ea \u05e8\u05d0\u05e0\u05e7 \u05e0\u05de\u05d5\u05da (\u05de\u05db\u05e4\u05dc\u05d4 \u05e9\u05dc \u05e9\u05ea\u05d9 \u05de\u05d8\u05e8\u05d9\u05e6\u05d5\u05ea \u05de\u05dc\u05d1\u05e0\u05d9\u05d5\u05ea \u05db\u05d0\u05e9\u05e8 \u05d4\u05de\u05d9\u05de\u05d3 \u05d4\u05e4\u05e0\u05d9\u05de\u05d9 \u05e9\u05dc \u05d4\u05de\u05db\u05e4\u05dc\u05d4 \u05e0\u05de\u05d5\u05da - \u05db\u05dc\u05d5\u05de\u05e8 (M x k * k x N) \u05db\u05d0\u05e9\u05e8 k \u05e7\u05d8\u05df \u05d4\u05e8\u05d1\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05de- M \u05d5- \u05de-M. \u05de\u05d7\u05e9\u05d1\u05d9\u05dd \u05d0\u05ea \u05d4\u05e7\u05dc\u05d8 \u05dc\u05e1\u05d5\u05e4\u05d8\u05de\u05e7\u05e1 \u05e2\u05dd \u05d4\u05de\u05d8\u05e8\u05d9\u05e6\u05d5\u05ea \u05d4\u05d0\u05dc\u05d5. \u05dc\u05d0\u05d7\u05e8 \u05de\u05db\u05df \u05de\u05e9\u05e8\u05e9\u05e8\u05d9\u05dd \u05d0\u05d5\u05ea\u05dd \u05e2\u05dd \u05d4\u05e7\u05dc\u05d8 \u05dc\u05e1\u05d5\u05e4\u05d8\u05de\u05e7\u05e1 \u05de\u05d4\u05e9\u05db\u05d1\u05d4 \u05d4\u05e7\u05d5\u05d3\u05de\u05ea, \u05de\u05e4\u05e2\u05d9\u05dc\u05d9\u05dd FFN \u05d5\u05d4\u05e0\u05d4 \u05d9\u05e9 \u05dc\u05e0\u05d5 \u05e7\u05dc\u05d8 \u05dc\u05e1\u05d5\u05e4\u05d8\u05de\u05e7\u05e1 \u05d1\u05e9\u05db\u05d1\u05d4 n. \u05d5\u05e9\u05d9\u05de\u05d5 \u05dc\u05d1 \u05e9\u05d0\u05e0\u05d5 \u05e6\u05e8\u05d9\u05db\u05d9\u05dd \u05dc\u05e9\u05de\u05d5\u05e8 \u05d4\u05e8\u05d1\u05d4 \u05e4\u05d7\u05d5\u05ea \u05d3\u05d0\u05d8\u05d4 \u05d1- KV-cache \u05db\u05d9 \u05d9\u05e9 \u05dc\u05e0\u05d5 \u05de\u05d8\u05e8\u05d9\u05e6\u05d5\u05ea \u05d1\u05e2\u05dc\u05d5\u05ea \u05e8\u05d0\u05e0\u05e7 \u05e0\u05de\u05d5\u05da.\", \"Replace\");\nparagraphs.items[6].insertText(\"\u05d0\u05d9\u05da \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05d0\u05ea \u05d4\u05e1\u05d9\u05e4\u05d5\u05e8 \u05d4\u05d6\u05d4? \u05de\u05e9\u05dc\u05d1\u05d9\u05dd \u05d0\u05ea \u05d4\u05dc\u05d5\u05e1 \u05d4\u05e8\u05d2\u05d9\u05dc \u05e9\u05dc \u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4 \u05e2\u05dd \u05dc\u05d5\u05e1 distillation \u05e9\u05de\u05d8\u05e8\u05ea\u05d4 \u05dc\u05e7\u05e8\u05d1 \u05d0\u05ea \u05de\u05e7\u05d3\u05de\u05d9 -attention \u05d4\u05de\u05d7\u05d5\u05e9\u05d1\u05d9\u05dd \u05d1\u05d3\u05e8\u05da \u05d4\u05de\u05d5\u05e6\u05e2\u05ea \u05e2\u05dd \u05d0\u05dc\u05d5 \u05e9\u05de\u05d7\u05d5\u05e9\u05d1\u05d9\u05dd \u05e2\u05dd \u05de\u05d5\u05d3\u05dc \u05e8\u05d2\u05d9\u05dc (\u05e2\u05dd attention \u05d5- KV Cache \u05e8\u05d2\u05d9\u05dc\u05d9\u05dd).\", \"Replace\");\n\n// Append two new paragraphs at the end of the body\nbody.insertParagraph(\"\u05de\u05d0\u05de\u05e8 \u05d3\u05d9 \u05de\u05e2\u05e0\u05d9\u05d9\u05df - \u05d0\u05d1\u05dc \u05e7\u05e6\u05ea \u05d0\u05e8\u05d5\u05da \u05de\u05d3\u05d9 \u05dc\u05d3\u05e2\u05ea\u05d9 \u05d0\u05d6 \u05ea\u05de\u05e6\u05ea\u05ea\u05d9 \u05dc\u05db\u05dd \u05d0\u05d5\u05ea\u05d5 \ud83d\ude42\", \"End\");\nbody.insertParagraph(\"https://arxiv.org/abs/2408.01890\", \"End\");\n\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n# Paragraph 1 (1-based): title/date line - replace whole range text\n# (collapses the existing run + <w:br/> + second run into a single run)\n$d.Paragraphs(1).Range.Text = '\u26a1\ufe0f\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 -07.11.24: \u26a1\ufe0f\ud83d\ude80'\n\n# Paragraphs 2-7 (1-based): replace full paragraph text\n$d.Paragraphs(2).Range.Text = 'Cross-layer Attention Sharing for Large Language Models'\n$d.Paragraphs(3).Range.Text = '\u05d0\u05ea\u05dd \u05d1\u05d8\u05d7 \u05d9\u05d5\u05d3\u05e2\u05d9\u05dd \u05d4\u05e8\u05e6\u05d4 \u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4 \u05e2\u05dc\u05d5\u05dc \u05dc\u05d4\u05d9\u05d5\u05ea \u05d3\u05d1\u05e8 \u05d3\u05d9 \u05d9\u05e7\u05e8 \u05de\u05d1\u05d7\u05d9\u05e0\u05ea \u05de\u05e9\u05d0\u05d1\u05d9 \u05d7\u05d9\u05e9\u05d5\u05d1 \u05d5\u05d2\u05dd \u05d4\u05d6\u05db\u05e8\u05d5\u05df. \u05d1\u05d8\u05d7 \u05db\u05d0\u05e9\u05e8 \u05d9\u05e9 \u05dc\u05db\u05dd \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05e2\u05dd \u05e2\u05e9\u05e8\u05d5\u05ea \u05de\u05d9\u05dc\u05d9\u05d0\u05e8\u05d3\u05d9 \u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd \u05e2\u05dc \u05e2\u05e9\u05e8\u05d5\u05ea \u05e8\u05d1\u05d5\u05ea \u05e9\u05dc \u05e9\u05db\u05d1\u05d5\u05ea \u05e9\u05dc \u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8\u05d9\u05dd. \u05d0\u05d7\u05d3 \u05d4\u05d3\u05d1\u05e8\u05d9\u05dd \u05d4\u05db\u05d1\u05d3\u05d9\u05dd \u05e9\u05de\u05e6\u05e8\u05d9\u05db\u05d9\u05dd \u05dc\u05d0 \u05de\u05e2\u05d8 \u05d6\u05d9\u05db\u05e8\u05d5\u05df \u05d4\u05d5\u05d0 KV-Cache, \u05e9\u05d1\u05d5 \u05e0\u05e9\u05de\u05e8\u05d9\u05dd \u05d4\u05de\u05db\u05e4\u05dc\u05d5\u05ea \u05e9\u05dc \u05d9\u05d9\u05e6\u05d5\u05d2\u05d9 (\u05d0\u05de\u05d1\u05d3\u05d9\u05e0\u05d2\u05e1) \u05e9\u05dc \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d1\u05de\u05d8\u05e8\u05d9\u05e6\u05d5\u05ea K \u05d5- V \u05dc\u05db\u05dc \u05d4\u05e9\u05db\u05d1\u05d5\u05ea \u05d5\u05dc\u05db\u05dc \u05d4\u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05e9\u05db\u05d1\u05e8 \u05d2\u05d5\u05e0\u05e8\u05d8\u05d5 (\u05db\u05d5\u05dc\u05dc \u05d4\u05e4\u05e8\u05d5\u05de\u05e4\u05d8 - \u05de\u05d3\u05d5\u05d1\u05e8 \u05d1\u05de\u05d5\u05d3\u05dc\u05d9 \u05d4\u05d3\u05e7\u05d5\u05d3\u05e8\u05d9\u05dd).'\n$d.Paragraphs(4).Range.Text = '\u05db\u05de\u05d5\u05d1\u05df \u05e9\u05db\u05d0\u05e9\u05e8 \u05d4\u05de\u05d9\u05de\u05d3\u05d9\u05dd \u05e9\u05dc \u05d5\u05e7\u05d8\u05d5\u05e8\u05d9 \u05d4\u05d9\u05d9\u05e6\u05d5\u05d2 \u05d5\u05d4\u05de\u05d8\u05e8\u05d9\u05e6\u05d5\u05ea \u05dc\u05d0 \u05e7\u05d8\u05e0\u05d9\u05dd \u05d5\u05d2\u05dd \u05d0\u05d5\u05e8\u05da \u05d4\u05d4\u05e7\u05e9\u05e8 \u05e0\u05de\u05d3\u05d3 \u05d1\u05e2\u05e9\u05e8\u05d5\u05ea \u05d5\u05de\u05d0\u05d5\u05ea \u05d0\u05dc\u05e4\u05d9\u05dd KV-Cache \u05d3\u05d5\u05e8\u05e9 \u05d4\u05e8\u05d1\u05d4 \u05de\u05d0\u05d5\u05d3 \u05d6\u05d9\u05db\u05e8\u05d5\u05df. \u05d1\u05e2\u05d1\u05e8 \u05d9\u05e6\u05d0\u05d5 \u05dc\u05d0 \u05de\u05e2\u05d8 \u05de\u05d0\u05de\u05e8\u05d9\u05dd \u05e9\u05e0\u05d9\u05e1\u05d5 \u05dc\u05d3\u05d7\u05d5\u05e1 \u05d0\u05d5\u05ea\u05d5 \u05e2\u05dc \u05d9\u05d3\u05d9 \u05e0\u05d9\u05ea\u05d5\u05d7 \u05d5\u05d6\u05d9\u05d4\u05d5\u05d9 \u05d9\u05ea\u05d9\u05e8\u05d5\u05d9\u05d5\u05ea \u05d0\u05d1\u05dc \u05d6\u05d4 \u05d1\u05d3\u05f4\u05db \u05e0\u05e2\u05e9\u05d4 \u05e4\u05e8 \u05e9\u05db\u05d1\u05d4 (= \u05d1\u05dc\u05d5\u05e7 \u05d4\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8). \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05de\u05e1\u05d5\u05e7\u05e8 \u05de\u05e6\u05d9\u05e2 \u05dc\u05d4\u05ea\u05d1\u05d5\u05e0\u05df \u05d1\u05d3\u05d7\u05d9\u05e1\u05ea KV-cache \u05de\u05e4\u05e8\u05e1\u05e4\u05e7\u05d8\u05d9\u05d1\u05d4 \u05e8\u05d7\u05d1\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05d5\u05dc\u05e0\u05e1\u05d5\u05ea \u05dc\u05d3\u05d7\u05d5\u05e1 \u05d0\u05d5\u05ea\u05d5 \u05d3\u05e8\u05da \u05e0\u05d9\u05e6\u05d5\u05dc \u05d4\u05ea\u05dc\u05d5\u05d9\u05d5\u05ea \u05e9\u05dc \u05d4-KV-cache \u05d1\u05d9\u05df \u05d4\u05e9\u05db\u05d1\u05d5\u05ea \u05d4\u05e9\u05d5\u05e0\u05d5\u05ea.'\n$d.Paragraphs(5).Range.Text = '\u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05d7\u05e7\u05e8\u05d5 \u05d3\u05de\u05d9\u05d5\u05df \u05d1\u05d9\u05df \u05d4\u05d7\u05dc\u05e7\u05d9\u05dd \u05d4\u05e9\u05d5\u05e0\u05d9\u05dd \u05d1\u05d1\u05dc\u05d5\u05e7 \u05d4\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8\u05d9\u05dd  (\u05de\u05db\u05e4\u05dc\u05d5\u05ea \u05e9\u05dc \u05d4\u05de\u05d8\u05e8\u05d9\u05e6\u05d5\u05ea \u05d4\u05e9\u05d5\u05e0\u05d5\u05ea \u05d1\u05d5\u05e7\u05d8\u05d5\u05e8\u05d9 \u05d9\u05d9\u05e6\u05d5\u05d2, \u05de\u05e7\u05d3\u05de\u05d9 attention \u05d5\u05db\u05d3\u05d5\u05de\u05d4) \u05d5\u05d4\u05d2\u05d9\u05e2\u05d5 \u05dc\u05de\u05e1\u05e7\u05e0\u05d4 \u05e9\u05e0\u05d9\u05ea\u05df \u05f4\u05dc\u05d4\u05e1\u05d9\u05e7\u05f4 \u05d0\u05ea \u05de\u05e7\u05d3\u05de\u05d9 \u05d4-attention \u05e9\u05dc \u05e9\u05db\u05d1\u05d4 n \u05de\u05d4\u05d3\u05d0\u05d8\u05d4 \u05e9\u05dc \u05e9\u05db\u05d1\u05d4 n-1 \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d7\u05e1\u05db\u05d5\u05e0\u05d9\u05ea \u05d7\u05d9\u05e9\u05d5\u05d1\u05d9\u05ea. \u05db\u05dc\u05d5\u05de\u05e8 \u05e2\u05dd \u05d4\u05e8\u05d1\u05d4 \u05e4\u05d7\u05d5\u05ea \u05de\u05e9\u05e7\u05d5\u05dc\u05d5\u05ea \u05de\u05d4\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05de\u05e8 \u05d4\u05e8\u05d2\u05d9\u05dc. \u05db\u05dc\u05d5\u05de\u05e8 \u05d4\u05d4\u05e6\u05e2\u05d4 \u05d4\u05d9\u05d0 \u05dc\u05e2\u05e9\u05d5\u05ea \u05e1\u05d5\u05d2 \u05e9\u05dc  LoRa \u05d0\u05d1\u05dc \u05dc\u05de\u05e7\u05d3\u05de\u05d9 \u05d4-attention. '\n$d.Paragraphs(6).Range.Text = '\u05d1\u05e6\u05d5\u05e8\u05d4 \u05e7\u05e6\u05ea \u05d9\u05d5\u05ea\u05e8 \u05e7\u05d5\u05e0\u05e7\u05e8\u05d8\u05d9\u05ea \u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d7\u05dc\u05d9\u05e3 \u05de\u05d8\u05e8\u05d9\u05e6\u05d5\u05ea W_Q \u05d5-W_K \u05d1\u05de\u05d8\u05e8\u05d9\u05e6\u05d5\u05ea \u05d1\u05e2\u05dc\u05d5\u05ea \u05e8\u05d0\u05e0\u05e7 \u05e0\u05de\u05d5\u05da (\u05de\u05db\u05e4\u05dc\u05d4 \u05e9\u05dc \u05e9\u05ea\u05d9 \u05de\u05d8\u05e8\u05d9\u05e6\u05d5\u05ea \u05de\u05dc\u05d1\u05e0\u05d9\u05d5\u05ea \u05db\u05d0\u05e9\u05e8 \u05d4\u05de\u05d9\u05de\u05d3 \u05d4\u05e4\u05e0\u05d9\u05de\u05d9 \u05e9\u05dc \u05d4\u05de\u05db\u05e4\u05dc\u05d4 \u05e0\u05de\u05d5\u05da - \u05db\u05dc\u05d5\u05de\u05e8 (M x k * k x N) \u05db\u05d0\u05e9\u05e8 k \u05e7\u05d8\u05df \u05d4\u05e8\u05d1\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05de- M \u05d5- \u05de-M. \u05de\u05d7\u05e9\u05d1\u05d9\u05dd \u05d0\u05ea \u05d4\u05e7\u05dc\u05d8 \u05dc\u05e1\u05d5\u05e4\u05d8\u05de\u05e7\u05e1 \u05e2\u05dd \u05d4\u05de\u05d8\u05e8\u05d9\u05e6\u05d5\u05ea \u05d4\u05d0\u05dc\u05d5. \u05dc\u05d0\u05d7\u05e8 \u05de\u05db\u05df \u05de\u05e9\u05e8\u05e9\u05e8\u05d9\u05dd \u05d0\u05d5\u05ea\u05dd \u05e2\u05dd \u05d4\u05e7\u05dc\u05d8 \u05dc\u05e1\u05d5\u05e4\u05d8\u05de\u05e7\u05e1 \u05de\u05d4\u05e9\u05db\u05d1\u05d4 \u05d4\u05e7\u05d5\u05d3\u05de\u05ea, \u05de\u05e4\u05e2\u05d9\u05dc\u05d9\u05dd FFN \u05d5\u05d4\u05e0\u05d4 \u05d9\u05e9 \u05dc\u05e0\u05d5 \u05e7\u05dc\u05d8 \u05dc\u05e1\u05d5\u05e4\u05d8\u05de\u05e7\u05e1 \u05d1\u05e9\u05db\u05d1\u05d4 n. \u05d5\u05e9\u05d9\u05de\u05d5 \u05dc\u05d1 \u05e9\u05d0\u05e0\u05d5 \u05e6\u05e8\u05d9\u05db\u05d9\u05dd \u05dc\u05e9\u05de\u05d5\u05e8 \u05d4\u05e8\u05d1\u05d4 \u05e4\u05d7\u05d5\u05ea \u05d3\u05d0\u05d8\u05d4 \u05d1- KV-cache \u05db\u05d9 \u05d9\u05e9 \u05dc\u05e0\u05d5 \u05de\u05d8\u05e8\u05d9\u05e6\u05d5\u05ea \u05d1\u05e2\u05dc\u05d5\u05ea \u05e8\u05d0\u05e0\u05e7 \u05e0\u05de\u05d5\u05da.'\n$d.Paragraphs(7).Range.Text = '\u05d0\u05d9\u05da \u05de\u05d0\u05de\u05e0\u05d9\u05dd \u05d0\u05ea \u05d4\u05e1\u05d9\u05e4\u05d5\u05e8 \u05d4\u05d6\u05d4? \u05de\u05e9\u05dc\u05d1\u05d9\u05dd \u05d0\u05ea \u05d4\u05dc\u05d5\u05e1 \u05d4\u05e8\u05d2\u05d9\u05dc \u05e9\u05dc \u05de\u05d5\u05d3\u05dc \u05e9\u05e4\u05d4 \u05e2\u05dd \u05dc\u05d5\u05e1 distillation \u05e9\u05de\u05d8\u05e8\u05ea\u05d4 \u05dc\u05e7\u05e8\u05d1 \u05d0\u05ea \u05de\u05e7\u05d3\u05de\u05d9 -attention \u05d4\u05de\u05d7\u05d5\u05e9\u05d1\u05d9\u05dd \u05d1\u05d3\u05e8\u05da \u05d4\u05de\u05d5\u05e6\u05e2\u05ea \u05e2\u05dd \u05d0\u05dc\u05d5 \u05e9\u05de\u05d7\u05d5\u05e9\u05d1\u05d9\u05dd \u05e2\u05dd \u05de\u05d5\u05d3\u05dc \u05e8\u05d2\u05d9\u05dc (\u05e2\u05dd attention \u05d5- KV Cache \u05e8\u05d2\u05d9\u05dc\u05d9\u05dd).'\n\n# Append two new paragraphs at the end of the document\n$d.Paragraphs($d.Paragraphs.Count).Range.InsertParagraphAfter()\n$d.Paragraphs($d.Paragraphs.Count).Range.Text = '\u05de\u05d0\u05de\u05e8 \u05d3\u05d9 \u05de\u05e2\u05e0\u05d9\u05d9\u05df - \u05d0\u05d1\u05dc \u05e7\u05e6\u05ea \u05d0\u05e8\u05d5\u05da \u05de\u05d3\u05d9 \u05dc\u05d3\u05e2\u05ea\u05d9 \u05d0\u05d6 \u05ea\u05de\u05e6\u05ea\u05ea\u05d9 \u05dc\u05db\u05dd \u05d0\u05d5\u05ea\u05d5 \ud83d\ude42'\n$d.Paragraphs($d.Paragraphs.Count).Range.InsertParagraphAfter()\n$d.Paragraphs($d.Paragraphs.Count).Range.Text = 'https://arxiv.org/abs/2408.01890'\n"}
